$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row update: insert a new "LB" header at D1, shift obj/CI/train_time/test_time right ---
$ws.Range("H1").Value = "test_time"
$ws.Range("G1").Value = "train_time"
$ws.Range("F1").Value = "CI"
$ws.Range("E1").Value = "obj"
$ws.Range("D1").Value = "LB"

# --- Data rows: update numeric results (in place, no column shift) ---
$ws.Range("D2").Value = 7234279.07865653
$ws.Range("E2").Value = 6342880.4636573996
$ws.Range("F2").Value = 873099.41119213204
$ws.Range("G2").Value = 43.1936388015747
$ws.Range("H2").Value = 5.06329274177551
$ws.Range("D3").Value = 3709850.2920581801
$ws.Range("E3").Value = 3455676.9370984202
$ws.Range("F3").Value = 212526.66464894699
$ws.Range("G3").Value = 131.73842382430999
$ws.Range("H3").Value = 6.9414899349212602
$ws.Range("D4").Value = 3279854.349891
$ws.Range("E4").Value = 3100418.8371999399
$ws.Range("F4").Value = 183713.694034418
$ws.Range("G4").Value = 99.673047065734806
$ws.Range("H4").Value = 6.2072639465331996
$ws.Range("D5").Value = 3140349.4965247302
$ws.Range("E5").Value = 2959960.3577384902
$ws.Range("F5").Value = 167100.84757933699
$ws.Range("G5").Value = 124.717630863189
$ws.Range("H5").Value = 6.4844388961791903
$ws.Range("D6").Value = 3090671.6440334101
$ws.Range("E6").Value = 2922870.4704522998
$ws.Range("F6").Value = 166407.21527551999
$ws.Range("G6").Value = 138.55235815048201
$ws.Range("H6").Value = 6.6143000125885001
$ws.Range("D7").Value = 4888378.3609059099
$ws.Range("E7").Value = 4832407.8515290804
$ws.Range("F7").Value = 527170.20557321305
$ws.Range("G7").Value = 24.474807024002001
$ws.Range("H7").Value = 5.4679908752441397
$ws.Range("D8").Value = 3421149.09784922
$ws.Range("E8").Value = 3322949.1143276799
$ws.Range("F8").Value = 225958.40904065399
$ws.Range("G8").Value = 41.296460151672299
$ws.Range("H8").Value = 6.2209599018096897
$ws.Range("D9").Value = 3124531.7976661902
$ws.Range("E9").Value = 3078117.6668686401
$ws.Range("F9").Value = 263259.74107469001
$ws.Range("G9").Value = 61.789297342300401
$ws.Range("H9").Value = 6.77239966392517
$ws.Range("D10").Value = 2975301.1958625098
$ws.Range("E10").Value = 2950489.6962703601
$ws.Range("F10").Value = 267698.29052467801
$ws.Range("G10").Value = 65.145031929016099
$ws.Range("H10").Value = 6.8974089622497496
$ws.Range("D11").Value = 2918795.8360245801
$ws.Range("E11").Value = 2893697.12872244
$ws.Range("F11").Value = 264198.90552844002
$ws.Range("G11").Value = 37.393470048904398
$ws.Range("H11").Value = 5.9280660152435303

# --- Selection change to match the saved view state ---
$ws.Range("D2").Select()
